# "Add files via upload" — a new data row (row 5) was appended to Sheet1,
# duplicating the submission that was already on row 2 (Submission ID
# "qJPar8" / Respondent ID "jBVv5Q" / andre.amorim@planejamento.mg.gov.br).
#
# Reproduce it by duplicating row 2 (values + formatting) down onto row 5.

$xlPasteFormats = -4122
$xlPasteValues  = -4163

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting first, then values, so the new row's number formats
# (date/plain-number columns) match row 2's without creating new style
# entries in the workbook.
$ws.Range("A2:V2").Copy()
$ws.Range("A5:V5").PasteSpecial($xlPasteFormats)

$ws.Range("A2:V2").Copy()
$ws.Range("A5:V5").PasteSpecial($xlPasteValues)

$ws.Range("A5:V5").Select()
